# "update data + download"
# After refreshing/extending the data in columns A, C and D, Excel's column
# widths were re-fit to the (now longer) content and the workbook was left
# scrolled down with a new active selection before being saved/downloaded
# again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (best-fit re-applied after the data refresh) ---------
# Column A (Class) widened slightly: 8.85546875 -> 9 characters.
$ws.Columns.Item(1).ColumnWidth = 8.14

# Column C (Specialized) now needs the same best-fit width as column B
# (Type): 11 characters.
$ws.Columns.Item(3).ColumnWidth = 10.15

# Column D (Detail-Specialized) grew to fit its longer values: ~17.29
# characters.
$ws.Columns.Item(4).ColumnWidth = 16.43

# --- View state: scrolled down and a new cell selected ------------------
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("K23").Select()
